$wb = $excel.ActiveWorkbook

# "4c55c43b-8590-44e4-9b76-677070aaad66.md" has been handed back and is now
# in sync with en-US. Update its status / timestamps / error detail across
# the Overview, zh-cn and de-de report sheets (row 3 in each case).

# Overview sheet: update status for the handed-back file (row 3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: update status, latest handback datetime, clear error detail
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("K3").Value = "2016-10-24 06:35:06"
$wsZh.Range("P3").Value = ""
$wsZh.Columns.Item(16).ColumnWidth = 12.86

# de-de sheet: update status, latest handback datetime, clear error detail
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("K3").Value = "2016-10-24 06:35:22"
$wsDe.Range("P3").Value = ""
$wsDe.Columns.Item(16).ColumnWidth = 12.86
